$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top of the sheet, shifting existing data down
$ws.Rows.Item(1).Insert()

# Populate the new header row with Key / Value labels
$ws.Cells.Item(1, 1).Value = "Key"
$ws.Cells.Item(1, 2).Value = "Value"

# Update the active selection to match the target state
$ws.Range("C10").Select()
